$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial number (45181) for every data
# row (2 through 410). Update it to 45182 for all of them.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 410) { $lastRow = 410 }

$range = $ws.Range("C2:C$lastRow")
$range.Value = 45182
